$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the coin price/volume table with the latest scraped values.
# Cells are forced to Text format while assigning so that price strings
# like "524.34" or "7.26" are stored verbatim instead of being
# reinterpreted by Excel as numbers; the style is then restored to
# "Normal" so no extra cell formatting is introduced.
function Set-TextValue($cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue 'D2' '57.668.51'
Set-TextValue 'E2' '  -0.09%  '
Set-TextValue 'D3' '3.118.73'
Set-TextValue 'E3' '  +1.35%  '
Set-TextValue 'E4' '  +0.01%  '
Set-TextValue 'D5' '524.34'
Set-TextValue 'E5' '  +1.04%  '
Set-TextValue 'D6' '141.02'
Set-TextValue 'E6' '  -1.20%  '
Set-TextValue 'E7' '  +0.08%  '
Set-TextValue 'D8' '3.120.44'
Set-TextValue 'E8' '  +1.49%  '
Set-TextValue 'D9' '0.433'
Set-TextValue 'E9' '  -0.29%  '
Set-TextValue 'D10' '7.26'
Set-TextValue 'E10' '  -0.05%  '
Set-TextValue 'E11' '  +1.33%  '
Set-TextValue 'D12' '0.386'
Set-TextValue 'E12' '  +2.51%  '
Set-TextValue 'D13' '3.659.05'
Set-TextValue 'E13' '  +1.59%  '
Set-TextValue 'E14' '  +1.68%  '
Set-TextValue 'D15' '26.32'
Set-TextValue 'E15' '  +2.41%  '
Set-TextValue 'E16' '  +0.87%  '
Set-TextValue 'D17' '57.760.18'
Set-TextValue 'E17' '  -0.02%  '
Set-TextValue 'D18' '3.120.94'
Set-TextValue 'E18' '  +1.52%  '
Set-TextValue 'E19' '  +0.29%  '
Set-TextValue 'D20' '12.87'
Set-TextValue 'E20' '  -0.13%  '
Set-TextValue 'D21' '8.07'
Set-TextValue 'E21' '  -0.86%  '
Set-TextValue 'D22' '337.07'
Set-TextValue 'E22' '  +1.32%  '
Set-TextValue 'E23' '  +0.00%  '
Set-TextValue 'D24' '0.511'
Set-TextValue 'E24' '  +2.37%  '
Set-TextValue 'D25' '66.78'
Set-TextValue 'E25' '  +1.37%  '
Set-TextValue 'E26' '  -0.34%  '
Set-TextValue 'E27' '  +0.12%  '
Set-TextValue 'D28' '0.0₃0924'
Set-TextValue 'E28' '  +2.20%  '
Set-TextValue 'E29' '  +3.49%  '
Set-TextValue 'D31' '7.21'
Set-TextValue 'E31' '  +0.02%  '
Set-TextValue 'D32' '1.87'
Set-TextValue 'E32' '  +2.84%  '
Set-TextValue 'D33' '1.20'
Set-TextValue 'E33' '  +0.06%  '
Set-TextValue 'D34' '20.94'
Set-TextValue 'E34' '  +0.77%  '
Set-TextValue 'D35' '155.71'
Set-TextValue 'E35' '  +0.71%  '
Set-TextValue 'E36' '  +3.29%  '
Set-TextValue 'D37' '6.12'
Set-TextValue 'E37' '  +2.99%  '
Set-TextValue 'D38' '27.12'
Set-TextValue 'E38' '  -0.04%  '
Set-TextValue 'E39' '  +1.61%  '
Set-TextValue 'D40' '0.0665'
Set-TextValue 'E40' '  -1.11%  '
Set-TextValue 'B41' 'Stacks'
Set-TextValue 'C41' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D41' '1.54'
Set-TextValue 'E41' '  +12.80%  '
Set-TextValue 'B42' 'RenzoRestakedETH'
Set-TextValue 'C42' 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
Set-TextValue 'D42' '3.164.34'
Set-TextValue 'E42' '  +1.58%  '
Set-TextValue 'D43' '0.691'
Set-TextValue 'E43' '  +5.45%  '
Set-TextValue 'D44' '3.92'
Set-TextValue 'E44' '  -0.51%  '
Set-TextValue 'D45' '36.87'
Set-TextValue 'E45' '  +0.36%  '
Set-TextValue 'E46' '  +0.00%  '
Set-TextValue 'D47' '2.293.17'
Set-TextValue 'E47' '  +1.15%  '
Set-TextValue 'D48' '0.0260'
Set-TextValue 'E48' '  +0.23%  '
Set-TextValue 'D49' '0.990'
Set-TextValue 'E49' '  +7.06%  '
Set-TextValue 'D50' '20.83'
Set-TextValue 'E50' '  -0.22%  '
Set-TextValue 'E51' '  +2.13%  '
